$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.232.59"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.438.58"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.70"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.80"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.92%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.435.91"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.63%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.97%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.72%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.24%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.24%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.127.48"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.437.14"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.88%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.65%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.18"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.84%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.96%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.60"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.24"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "629.16"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.559.19"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "PEPE"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0948"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -9.87%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Binance-PegBSC-USD"

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.01%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.18%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.43"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.47%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.10"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.29%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.22"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.08%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.75"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.71"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.46%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0521"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.75%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.55"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -9.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0233"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.66%  "

